# "Small fix on formatted agreement"
#
# The template has a set of {{merge_field}} placeholders. One of them -
# {{tenant_addresss}} - lost its closing "}}" (it is split across runs
# because of the spell-checker's spellStart/spellEnd wrapping around the
# field name). This adds a new run containing the missing "}}" right
# after the field-name run, matching the Times New Roman / bCs run
# formatting already used by its sibling runs in that placeholder.

$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("tenant_addresss", $true, $false, $false, $false, `
                            $false, $true, 1, $false, "", 0)

if ($found) {
    # Collapse the found range to its end (right after "tenant_addresss",
    # i.e. right after the spellStart run, before spellEnd/</w:p>) and
    # insert the missing closing braces as a brand new run there.
    $rng.Collapse(0)   # wdCollapseEnd
    $rng.InsertAfter("}}")

    # Match the surrounding placeholder runs' formatting exactly:
    # <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman"
    #           w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:bCs/>
    $rng.Font.NameAscii = "Times New Roman"
    $rng.Font.NameFarEast = "Times New Roman"
    $rng.Font.NameOther = "Times New Roman"
    $rng.Font.NameBi = "Times New Roman"
    $rng.Font.BoldBi = $true
}

# styles.xml: the "Normal Table" built-in style picks up an explicit
# <w:unhideWhenUsed/> on this revision (part of a broader semiHidden /
# unhideWhenUsed normalization of the built-in styles). Apply the part
# of that normalization that is reachable through the Style object.
$tableNormal = $d.Styles("Normal Table")
$tableNormal.UnhideWhenUsed = $true
